$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final target table for rows 2..17 (A=index, B=name, C=from_bus, D=to_bus, E=in_service).
# Two new lines (line7, line8) were inserted before the extr* rows, and the
# from_bus/to_bus/in_service figures for several rows were recomputed.
$data = @(
  @(0,  "line1", 7,  9,  $true),
  @(1,  "line2", 9,  8,  $true),
  @(2,  "line3", 8,  10, $true),
  @(3,  "line4", 8,  11, $false),
  @(4,  "line5", 10, 5,  $true),
  @(5,  "line6", 12, 8,  $true),
  @(6,  "line7", 14, 11, $true),
  @(7,  "line8", 16, 9,  $true),
  @(8,  "extr1", 5,  12, $true),
  @(9,  "extr2", 5,  9,  $true),
  @(10, "extr3", 10, 11, $false),
  @(11, "extr4", 7,  8,  $true),
  @(12, "extr5", 9,  11, $false),
  @(13, "extr6", 7,  11, $true),
  @(14, "extr7", 5,  7,  $false),
  @(15, "extr8", 8,  5,  $true)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $entry = $data[$i]

    # Rows 16-17 are brand new; give column A the same bold/border/centered
    # style already used by the rest of the index column (A2:A15) by copying
    # the formatting down from the row above before writing the new value.
    if ($row -gt 15) {
        $ws.Range("A" + ($row - 1)).Copy() | Out-Null
        $ws.Range("A" + $row).PasteSpecial(-4122) | Out-Null
    }

    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
}
